$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Update values on sheet "results"
# ---------------------------------------------------------------
$wsResults = $wb.Worksheets.Item("results")

$resultsData = @(
    @(0.5925925925925926, 0.6, 0.5555555555555556, 0.576923076923077),
    @(0.6296296296296297, 0.6521739130434783, 0.5555555555555556, 0.6),
    @(0.5185185185185185, 0.5238095238095238, 0.4074074074074074, 0.4583333333333333),
    @(0.4444444444444444, 0.4482758620689655, 0.4814814814814815, 0.4642857142857143),
    @(0.4444444444444444, 0.4516129032258064, 0.5185185185185185, 0.4827586206896551),
    @(0.462962962962963,  0.46875,             0.5555555555555556, 0.5084745762711864),
    @(0.4444444444444444, 0.4285714285714285, 0.3333333333333333, 0.375)
)

for ($i = 0; $i -lt $resultsData.Length; $i++) {
    $row = $i + 2
    $vals = $resultsData[$i]
    $wsResults.Cells.Item($row, 2).Value = $vals[0]
    $wsResults.Cells.Item($row, 3).Value = $vals[1]
    $wsResults.Cells.Item($row, 4).Value = $vals[2]
    $wsResults.Cells.Item($row, 5).Value = $vals[3]
}

# ---------------------------------------------------------------
# 2) Update values on sheet "results with ST"
# ---------------------------------------------------------------
$wsST = $wb.Worksheets.Item("results with ST")

$stData = @(
    @(0.6296296296296297, 0.6296296296296297, 0.6296296296296297, 0.6296296296296297),
    @(0.5555555555555556, 0.5517241379310345, 0.5925925925925926, 0.5714285714285714),
    @(0.6666666666666666, 0.6666666666666666, 0.6666666666666666, 0.6666666666666666),
    @(0.5185185185185185, 0.52,                0.4814814814814815, 0.5),
    @(0.4259259259259259, 0.4230769230769231, 0.4074074074074074, 0.4150943396226415),
    @(0.462962962962963,  0.4666666666666667, 0.5185185185185185, 0.4912280701754386),
    @(0.5, 0.5, 0.7037037037037037, 0.5846153846153846)
)

for ($i = 0; $i -lt $stData.Length; $i++) {
    $row = $i + 2
    $vals = $stData[$i]
    $wsST.Cells.Item($row, 2).Value = $vals[0]
    $wsST.Cells.Item($row, 3).Value = $vals[1]
    $wsST.Cells.Item($row, 4).Value = $vals[2]
    $wsST.Cells.Item($row, 5).Value = $vals[3]
}

# ---------------------------------------------------------------
# 3) Rename the existing empty "Sheet" to "results gray" and fill
#    it with the new gray-images test results.
# ---------------------------------------------------------------
$wsGray = $wb.Worksheets.Item("Sheet")
$wsGray.Name = "results gray"

$wsGray.Cells.Item(1, 1).Value = "backbone"
$wsGray.Cells.Item(1, 2).Value = "accuracy"
$wsGray.Cells.Item(1, 3).Value = "precision"
$wsGray.Cells.Item(1, 4).Value = "recall"
$wsGray.Cells.Item(1, 5).Value = "f1"

$grayNames = @("VGG-Face", "Facenet", "Facenet512", "OpenFace", "DeepFace", "DeepID", "ArcFace")

$grayData = @(
    @(0.6851851851851852, 0.6470588235294118, 0.8148148148148148, 0.7213114754098361),
    @(0.5925925925925926, 0.6,                 0.5555555555555556, 0.576923076923077),
    @(0.6111111111111112, 0.7142857142857143, 0.3703703703703703, 0.4878048780487805),
    @(0.5555555555555556, 0.5384615384615384, 0.7777777777777778, 0.6363636363636364),
    @(0.4259259259259259, 0.4375,              0.5185185185185185, 0.4745762711864406),
    @(0.4814814814814815, 0.4871794871794872, 0.7037037037037037, 0.5757575757575758),
    @(0.4814814814814815, 0.4761904761904762, 0.3703703703703703, 0.4166666666666667)
)

for ($i = 0; $i -lt $grayData.Length; $i++) {
    $row = $i + 2
    $vals = $grayData[$i]
    $wsGray.Cells.Item($row, 1).Value = $grayNames[$i]
    $wsGray.Cells.Item($row, 2).Value = $vals[0]
    $wsGray.Cells.Item($row, 3).Value = $vals[1]
    $wsGray.Cells.Item($row, 4).Value = $vals[2]
    $wsGray.Cells.Item($row, 5).Value = $vals[3]
}

# ---------------------------------------------------------------
# 4) Add a new empty worksheet named "Sheet" after "results gray"
# ---------------------------------------------------------------
$wsNew = $wb.Worksheets.Add()
$wsNew.Name = "Sheet"
$wsGrayRef = $wb.Worksheets.Item("results gray")
$wsNew.Move($null, $wsGrayRef)
